$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the weekly remaining-work values for the tasks in rows 4-9 (column D)
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0.5
$ws.Range("D9").Value = 0

# Move the active cell selection from D4 to D3
$ws.Range("D3").Select()
